$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Product Name (title) cell to the new revision text
$ws.Range("C1").Value = "Mounting Bracket A 45 degree  (VER 1.1 REV A)"

# Update the BOM line item's Tag to the new version tag
$ws.Range("E7").Value = "MOUNTA.45.v1.1"

# Update "Last Updated" date (row 3) and the BOM line item's Date (row 7)
$ws.Range("F3").Value = 42397
$ws.Range("F7").Value = 42397

# Row 1 grew taller to accommodate the longer title text
$ws.Rows.Item(1).RowHeight = 55.5

# Move the saved selection to C2
$ws.Range("C2").Select()
